$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'310.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.11%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'48.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'8.37%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.245"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'2.08%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07808"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-3.01%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.533"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.19%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.315"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'21.24%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.563"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-7.03%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1251"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-3.42%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1950"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.91%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09360"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.29%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04542"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'7.22%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'0.31%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001307"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.48%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.04205"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.21%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005817"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.48%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'HotbitToken"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'0.004119"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-9.67%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "'LEO"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'3.336"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-1.66%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'BTSEToken"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'2.425"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.61%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'BitpandaEcosystemToken"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'0.3443"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'2.04%"
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "'MCDex"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'8.148"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-1.44%"
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "'ProBitToken"
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").Value = "'https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").Value = "'0.1363"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-0.64%"
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "'ZBToken"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'0.3069"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-2.25%"
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "'BitKan"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'0.001300"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.92%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001367"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.84%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003560"
$ws.Range("D26").Style = "Normal"
$ws.Range("D38").Value = "'0.02594"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-4.25%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05754"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'6.24%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.01035"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'88.74%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.008012"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.36%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1435"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.27%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.008409"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'14.65%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008662"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'0.85%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3112"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-0.46%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006953"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'2.35%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'1.13%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.05544"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-9.94%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.004024"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'1.17%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002112"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'1.13%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002011"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'1.13%"
$ws.Range("E51").Style = "Normal"
